# Weekly update: insert the newest week's observation for
# "Vega Central Mapocho de Santiago - Rabanito" as a new row 258,
# pushing the previously-top rows (258-272) down to (259-273).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row before the current row 258 (mirrors the weekly
# pattern already present in the sheet: newest entry lands on top of
# this block, everything else shifts down by one row).
$ws.Rows("258:258").Insert()

# Populate the new row with the latest week's data.
$ws.Range("A258").Value = 9
$ws.Range("B258").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C258").Value = "Metropolitana"
$ws.Range("D258").Value = 44714
$ws.Range("E258").Value = 13
$ws.Range("F258").Value = 300000001
$ws.Range("G258").Value = "Rabanito"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 5200
$ws.Range("K258").Value = 3000
$ws.Range("L258").Value = 3000
$ws.Range("M258").Value = 3000
$ws.Range("N258").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O258").Value = "Provincia de Chacabuco"
$ws.Range("P258").Value = 30
$ws.Range("Q258").Value = 100
$ws.Range("R258").Value = "Hortaliza"
